$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 164, pushing the existing rows 164-257
# down to 166-259 (matching the diff, which shows all rows from 166 onward
# taking on the values that used to belong to the row two positions above).
$ws.Rows.Item(164).Insert()
$ws.Rows.Item(164).Insert()

# Populate the first new row (164) with the new price-report entry.
$ws.Range("A164").Value = 3
$ws.Range("B164").Value = "Femacal de La Calera"
$ws.Range("C164").Value = "Coquimbo"
$ws.Range("D164").Value = 44518
$ws.Range("E164").Value = 5
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100108
$ws.Range("H164").Value = "Tropicales y subtropicales"
$ws.Range("I164").Value = 100108002
$ws.Range("J164").Value = "Mango"
$ws.Range("K164").Value = "Sin especificar"
$ws.Range("L164").Value = "Primera"
$ws.Range("M164").Value = 456
$ws.Range("N164").Value = 6500
$ws.Range("O164").Value = 6500
$ws.Range("P164").Value = 6500
$ws.Range("Q164").Value = '$/bandeja 4 kilos'
$ws.Range("R164").Value = "Perú"
$ws.Range("S164").Value = 1625
$ws.Range("T164").Value = 4

# Populate the second new row (165) with the new price-report entry.
$ws.Range("A165").Value = 3
$ws.Range("B165").Value = "Femacal de La Calera"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 44518
$ws.Range("E165").Value = 5
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100108
$ws.Range("H165").Value = "Tropicales y subtropicales"
$ws.Range("I165").Value = 100108002
$ws.Range("J165").Value = "Mango"
$ws.Range("K165").Value = "Sin especificar"
$ws.Range("L165").Value = "Segunda"
$ws.Range("M165").Value = 228
$ws.Range("N165").Value = 6500
$ws.Range("O165").Value = 6500
$ws.Range("P165").Value = 6500
$ws.Range("Q165").Value = '$/bandeja 4 kilos'
$ws.Range("R165").Value = "Perú"
$ws.Range("S165").Value = 1625
$ws.Range("T165").Value = 4
